$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.921.44'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '1.993.57'
$ws.Range('E3').Value = '  -2.88%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.59'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.640'
$ws.Range('E6').Value = '  -3.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.85'
$ws.Range('E7').Value = '  +9.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.65'
$ws.Range('E9').Value = '  -3.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.365'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0741'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('E13').Value = '  -2.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.78'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '2.283.98'
$ws.Range('E15').Value = '  -2.98%  '
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.37'
$ws.Range('E17').Value = '  +13.45%  '
$ws.Range('D18').Value = '1.990.23'
$ws.Range('E18').Value = '  -2.98%  '
$ws.Range('D19').Value = '35.866.65'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.86'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.22'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.54'
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.64'
$ws.Range('E25').Value = '  +16.49%  '
$ws.Range('E26').Value = '  -4.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.61'
$ws.Range('E27').Value = '  +4.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.11'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.41'
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.119'
$ws.Range('E30').Value = '  -1.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.93'
$ws.Range('E31').Value = '  -2.57%  '
$ws.Range('E32').Value = '  -5.54%  '
$ws.Range('E33').Value = '  +13.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0605'
$ws.Range('E34').Value = '  +1.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  +10.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.42'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.76'
$ws.Range('E39').Value = '  +14.62%  '
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0932'
$ws.Range('E43').Value = '  +1.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.10'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.57'
$ws.Range('E45').Value = '  +4.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '94.12'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.81'
$ws.Range('E47').Value = '  +3.85%  '
$ws.Range('D48').Value = '1.369.22'
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.33'
$ws.Range('E50').Value = '  +2.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.94'
$ws.Range('E51').Value = '  +2.56%  '
